$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update existing row 39 values (H39, I39)
$ws.Cells.Item(39, 8).Value = 29
$ws.Cells.Item(39, 9).Value = 156

# Add new row 40 with data from the press conference update
$ws.Cells.Item(40, 1).Value = 43927
$ws.Cells.Item(40, 1).NumberFormat = $ws.Cells.Item(39, 1).NumberFormat

$ws.Cells.Item(40, 2).Value = 39
$ws.Cells.Item(40, 3).Value = 911
$ws.Cells.Item(40, 4).Value = 28
$ws.Cells.Item(40, 5).Value = 195
$ws.Cells.Item(40, 6).Value = 67
$ws.Cells.Item(40, 7).Value = 1106
$ws.Cells.Item(40, 8).Value = 20
$ws.Cells.Item(40, 9).Value = 176
$ws.Cells.Item(40, 10).Value = 13
$ws.Cells.Item(40, 12).Value = 3
$ws.Cells.Item(40, 13).Value = 0
$ws.Cells.Item(40, 14).Value = 1
$ws.Cells.Item(40, 15).Value = 476
$ws.Cells.Item(40, 16).Value = 420
$ws.Cells.Item(40, 17).Value = 188
$ws.Cells.Item(40, 18).Value = 22
$ws.Cells.Item(40, 19).Value = 1106
$ws.Cells.Item(40, 20).Value = "Manual"
